# Insert a new data row before current row 128 (shifts existing rows 128..184
# down to 129..185) and populate it with the new record, replicating the
# other (unchanged) fields from the row that used to be at 128.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(128).Insert()

$ws.Cells.Item(128, 1).Value = 10
$ws.Cells.Item(128, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(128, 3).Value = "La Araucanía"
$ws.Cells.Item(128, 4).Value = 44609
$ws.Cells.Item(128, 5).Value = 9
$ws.Cells.Item(128, 6).Value = 100112043
$ws.Cells.Item(128, 7).Value = "Pepino dulce"
$ws.Cells.Item(128, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(128, 9).Value = "Primera"
$ws.Cells.Item(128, 10).Value = 150
$ws.Cells.Item(128, 11).Value = 23000
$ws.Cells.Item(128, 12).Value = 23000
$ws.Cells.Item(128, 13).Value = 23000
$ws.Cells.Item(128, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(128, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(128, 16).Value = 1278
$ws.Cells.Item(128, 17).Value = 18
$ws.Cells.Item(128, 18).Value = "Hortaliza"
